# Update the handback status report with freshly generated handoff/handback
# timestamps for the most recently processed file in each locale sheet.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 (9bfdd7b6-... file) gets new Correspond Handoff/Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 04:23:09"
$wsZhCn.Range("H3").Value = "2016-03-18 04:23:28"

# de-de sheet: row 3 (9bfdd7b6-... file) gets new Correspond Handoff/Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 04:23:12"
$wsDeDe.Range("H3").Value = "2016-03-18 04:23:32"
